$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly report row (row 10)
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2018.07.23"
$ws.Range("B10").Value = "建立接口，传入数据库"
$ws.Range("C10").Value = "建立接口，让输入的文本进行分析后跑进数据库中"
$ws.Range("D10").Value = "如何将网页输入的文本传回编译器中"

# A10 should match the same formatting as the A column "date" cells above it (horizontal+vertical center)
$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("A10").VerticalAlignment = -4108

# Leave selection on C10, matching the author's final cursor position
$ws.Range("C10").Select()
